$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet keeps a sliding-window log of "last updated" timestamps in
# column D, organized in three 14-row blocks (rows 2-15, 16-29, 30-43).
# A new "Actualizar" run pushes a fresh timestamp into the top block and
# shifts the two older blocks down by one slot; the oldest block falls off
# the bottom of the sheet.

$newest = 44239.63938445505
$middle = 44234.72344277778
$oldest = 44234.70222912037

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newest
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $middle
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldest
}
